$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3-4: re-categorize existing entries from "SEO" to "Accessibilité" ---
$ws.Range("A3").Value = "Accessibilité"
$ws.Range("A4").Value = "Accessibilité"

# --- Row 5: new Accessibilité entry (image compression) ---
$ws.Range("A5").Value = "Accessibilité"
$ws.Range("B5").Value = "Compression image"
$ws.Range("C5").Value = "Vérifier son code sur GTMetrix"
$ws.Range("D5").Value = "Résuolution d'image trop élévée donc ralentissement du chargement du site"
$ws.Range("E5").Value = "Vérifier que la taille de l'image initial soit adapter a l'image souhaité"

# --- Row 6: new Accessibilité entry (JS code) ---
$ws.Range("A6").Value = "Accessibilité"
$ws.Range("B6").Value = "Code JS"
$ws.Range("D6").Value = "Les espâces dans le code JS ralenti l'execution du code"
$ws.Range("E6").Value = "Compresser (minify) son code avec GTMetrix"

# --- Row 7: new Accessibilité entry (color contrast) ---
$ws.Range("A7").Value = "Accessibilité"
$ws.Range("B7").Value = "Contraste de couleur"
$ws.Range("C7").Value = "Pensez a tous les utilisateurs lors de la conception"
$ws.Range("D7").Value = "Tres faible contraste entre le texte et son background"
$ws.Range("E7").Value = "Appliquer les standards de contraste"

# --- Row 8: new Accessibilité entry (link problems) ---
$ws.Range("A8").Value = "Accessibilité"
$ws.Range("B8").Value = "Probleme de lien"
$ws.Range("C8").Value = "Verifier son code sur Wave"
$ws.Range("D8").Value = "Lien sans texte/Lien avec destination pas coherente..."
$ws.Range("E8").Value = "Directement mettre en place le lien de destination"
$ws.Rows.Item(8).RowHeight = 14

# --- Row 9: new Accessibilité entry (visual / responsive) ---
$ws.Range("A9").Value = "Accessibilité"
$ws.Range("B9").Value = "Visuel"
$ws.Range("C9").Value = "Tester son site sur differents appareils"
$ws.Range("D9").Value = "Plusieurs probleme d'affichage sur differentes tailles d'écran"
$ws.Range("E9").Value = "Bien travailler son responsive"

# --- Row 10: new SEO & Accessibilité entry (html tags) ---
$ws.Range("A10").Value = "SEO & Accessibilité"
$ws.Range("B10").Value = "Mauvaise balise html"
$ws.Range("C10").Value = "Bien identifer les parties du site avant de le construire"
$ws.Range("D10").Value = "Il n'y a que des balises <div> dans le code html"
$ws.Range("E10").Value = "Faire une maquette du site et identifier les balises"

# --- Row 11: new SEO entry (sitemap) ---
$ws.Range("A11").Value = "SEO"
$ws.Range("B11").Value = "Pas de fichier sitemap"
$ws.Range("C11").Value = "Vérifier son code sur des outils SEO"
$ws.Range("D11").Value = "Pas de fichier sitemap.xml pour optimiser l'interaction avec les bots"
$ws.Range("E11").Value = "Mettre un fichier sitemap.xml"
$ws.Range("F11").Value = "neilpatel.com"

# --- Row 12: trailing SEO category marker ---
$ws.Range("A12").Value = "SEO"

# --- Column width adjustments ---
# (Target OOXML raw widths are 40.765625 / 53.53515625; the COM
# ColumnWidth property quantizes to whole pixels at MDW=7, so these
# inputs are chosen to land on the closest achievable raw width.)
$ws.Columns.Item(3).ColumnWidth = 40
$ws.Columns.Item(4).ColumnWidth = 52.857142857142854

# --- Update the active selection to reflect last edited cell ---
$ws.Range("C23").Select()
